$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 293, shifting existing rows 293-345 down to 295-347
$ws.Rows("293:294").Insert()

# Populate the two newly inserted rows with new data (rows 293 and 294)
$ws.Range("A293").Value = 6
$ws.Range("B293").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C293").Value = 'Metropolitana'
$ws.Range("D293").Value = 44889
$ws.Range("E293").Value = 13
$ws.Range("F293").Value = 100112026
$ws.Range("G293").Value = 'Haba'
$ws.Range("H293").Value = 'Sin especificar'
$ws.Range("I293").Value = 'Primera'
$ws.Range("J293").Value = 580
$ws.Range("K293").Value = 7000
$ws.Range("L293").Value = 8000
$ws.Range("M293").Value = 7448
$ws.Range("N293").Value = '$/saco 25 kilos'
$ws.Range("O293").Value = 'Región de O''Higgins'
$ws.Range("P293").Value = 298
$ws.Range("Q293").Value = 25
$ws.Range("R293").Value = 'Hortaliza'
$ws.Range("A294").Value = 6
$ws.Range("B294").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C294").Value = 'Metropolitana'
$ws.Range("D294").Value = 44889
$ws.Range("E294").Value = 13
$ws.Range("F294").Value = 100112026
$ws.Range("G294").Value = 'Haba'
$ws.Range("H294").Value = 'Sin especificar'
$ws.Range("I294").Value = 'Primera'
$ws.Range("J294").Value = 500
$ws.Range("K294").Value = 6000
$ws.Range("L294").Value = 7000
$ws.Range("M294").Value = 6480
$ws.Range("N294").Value = '$/saco 25 kilos'
$ws.Range("O294").Value = 'Región del Maule'
$ws.Range("P294").Value = 259
$ws.Range("Q294").Value = 25
$ws.Range("R294").Value = 'Hortaliza'
